$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Give the (currently plain) leading run of six rich-text cells the same
#    explicit font formatting (11pt, black, 宋体) that the rest of their runs
#    already carry.  Excel stores these as shared strings, so Characters()
#    applied to the owning cell edits the rPr of that leading run only.
# ---------------------------------------------------------------------------
$richTextEdits = @(
    @{ Cell = "J7";  Length = 51 },
    @{ Cell = "F8";  Length = 20 },
    @{ Cell = "F9";  Length = 26 },
    @{ Cell = "J10"; Length = 24 },
    @{ Cell = "J11"; Length = 24 },
    @{ Cell = "F12"; Length = 30 }
)

foreach ($edit in $richTextEdits) {
    $chars = $ws.Range($edit.Cell).Characters(1, $edit.Length)
    $font = $chars.Font
    $font.Size = 11
    $font.Color = 0
    $font.Name = "宋体"
}

# ---------------------------------------------------------------------------
# 2) Window / view changes: zoom 70% -> 85%, selection moves to F6.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 85
$ws.Range("F6").Select()

# ---------------------------------------------------------------------------
# 3) New explicit width for column C.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 14.65

# ---------------------------------------------------------------------------
# 4) Row height adjustments.
# ---------------------------------------------------------------------------
$ws.Rows.Item(7).RowHeight = 86
$ws.Rows.Item(8).RowHeight = 58
$ws.Rows.Item(10).RowHeight = 40
$ws.Rows.Item(11).RowHeight = 65
$ws.Rows.Item(12).RowHeight = 51
$ws.Rows.Item(14).RowHeight = 14.75
$ws.Rows.Item(16).RowHeight = 55
